$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting: give E2 a new bold 10pt style (distinct from the other
# body cells) before changing its text, so the engine derives a clean new
# font/cellXf pair instead of inheriting the old wrap-text alignment block.
$ws.Range("E2").Style = "Normal"
$ws.Range("E2").Font.Bold = $true
$ws.Range("E2").Font.Size = 10

# --- Session content updates (slot grid reshuffle) ---
# Row 2 (Slot 1)
$ws.Range("E2").Value = "Chew: Pre-Skilling for the Unknown: Building a Human-Centered AI Culture in Schools"
$ws.Range("F2").Value = "Nurenberg: Dangerous (Artificial) Minds: Engaging student critical thinking and analysis in their interactions with AI"
$ws.Range("M2").Value = "Alsamadisi: Rethinking Creative and Critical Thinking in the Age of AI"

# Row 3 (Slot 2)
$ws.Range("L3").Value = "Place Based TBD"

# Row 4 (Slot 3)
$ws.Range("M4").Value = "MacClintic: Student Presentations "

# --- Selection moves to D3 ---
$ws.Range("D3").Select()
